$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$brk = [char]11

$cell = $t.Cell(1, 1)
$cell.Range.Text = [string]::Join($brk, @("54 x 69", "  6    9", "  ----", "5|    |", "4|    |"))
$cell = $t.Cell(1, 2)
$cell.Range.Text = [string]::Join($brk, @("90 x 96", "  9    6", "  ----", "9|    |", "0|    |"))
$cell = $t.Cell(1, 3)
$cell.Range.Text = [string]::Join($brk, @("22 x 24", "  2    4", "  ----", "2|    |", "2|    |"))
$cell = $t.Cell(2, 1)
$cell.Range.Text = [string]::Join($brk, @("65 x 90", "  9    0", "  ----", "6|    |", "5|    |"))
$cell = $t.Cell(2, 2)
$cell.Range.Text = [string]::Join($brk, @("59 x 20", "  2    0", "  ----", "5|    |", "9|    |"))
$cell = $t.Cell(2, 3)
$cell.Range.Text = [string]::Join($brk, @("47 x 36", "  3    6", "  ----", "4|    |", "7|    |"))
$cell = $t.Cell(3, 1)
$cell.Range.Text = [string]::Join($brk, @("59 x 36", "  3    6", "  ----", "5|    |", "9|    |"))
$cell = $t.Cell(3, 2)
$cell.Range.Text = [string]::Join($brk, @("51 x 41", "  4    1", "  ----", "5|    |", "1|    |"))
$cell = $t.Cell(3, 3)
$cell.Range.Text = [string]::Join($brk, @("36 x 85", "  8    5", "  ----", "3|    |", "6|    |"))
$cell = $t.Cell(4, 1)
$cell.Range.Text = [string]::Join($brk, @("81 x 46", "  4    6", "  ----", "8|    |", "1|    |"))
$cell = $t.Cell(4, 2)
$cell.Range.Text = [string]::Join($brk, @("36 x 51", "  5    1", "  ----", "3|    |", "6|    |"))
$cell = $t.Cell(4, 3)
$cell.Range.Text = [string]::Join($brk, @("27 x 76", "  7    6", "  ----", "2|    |", "7|    |"))
$cell = $t.Cell(5, 1)
$cell.Range.Text = [string]::Join($brk, @("32 x 78", "  7    8", "  ----", "3|    |", "2|    |"))
$cell = $t.Cell(5, 2)
$cell.Range.Text = [string]::Join($brk, @("48 x 20", "  2    0", "  ----", "4|    |", "8|    |"))
$cell = $t.Cell(5, 3)
$cell.Range.Text = [string]::Join($brk, @("40 x 11", "  1    1", "  ----", "4|    |", "0|    |"))
